$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    $ws.Cells.Item($row, 3).Value = 45180
}
